$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Wipe the old "Points / Nom de la BDF / Classement / ..." table and
# replace it with the single-column "Joueur" roster.
$ws.Range("A1:E2").ClearContents()

$values = @("Joueur","Baptiste","Eric","Côme","Didier","Sylvie P","Béa","Hugo D")
for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}

# Drop the leftover page-footer text from the old sheet.
$ps = $ws.PageSetup
$ps.LeftHeader = ""
$ps.CenterHeader = ""
$ps.RightHeader = ""
$ps.LeftFooter = ""
$ps.CenterFooter = ""
$ps.RightFooter = ""
$ps.OddHeader = ""
$ps.OddFooter = ""
$ps.EvenHeader = ""
$ps.EvenFooter = ""
$ps.FirstHeader = ""
$ps.FirstFooter = ""

# Restore the classic "Office" theme accent colors (accent1 <-> accent5
# were swapped when the sheet picked up the "Office 2013-2022" theme).
$colorScheme = $wb.Theme.ThemeColorScheme
$colorScheme.Colors(5).RGB = 13998939   # accent1 -> 5B9BD5
$colorScheme.Colors(9).RGB = 12874308   # accent5 -> 4472C4

# Move the active selection to A9, just past the new data.
$ws.Range("A9").Select() | Out-Null
